$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest six years (2004年-2009年), which sat in rows 2:7.
# Everything below shifts up six rows, so the former 2010年 row (8)
# becomes row 2, ... former 2020年 row (18) becomes row 12.
$ws.Rows("2:7").Delete()

# Append the new 2021年 row as row 13. Copy the previous last row (12,
# 2020年) first so the new row inherits matching cell styles/format,
# then overwrite with the 2021年 figures.
$ws.Range("A12:H12").Copy($ws.Range("A13:H13"))

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 3352364.344
$ws.Range("D13").Value = 728093.873
$ws.Range("E13").Value = 5436.764
$ws.Range("F13").Value = 36610.88
$ws.Range("G13").Value = 2220981.739
$ws.Range("H13").Value = 361241.088
